$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'245.79"
$ws.Range("E2").Value = "'-0.29%"
$ws.Range("D3").Value = "'30.12"
$ws.Range("E3").Value = "'-0.57%"
$ws.Range("D4").Value = "'5.156"
$ws.Range("E4").Value = "'-0.38%"
$ws.Range("D5").Value = "'0.05763"
$ws.Range("E5").Value = "'0.60%"
$ws.Range("E6").Value = "'1.02%"
$ws.Range("D7").Value = "'3.281"
$ws.Range("E7").Value = "'6.85%"
$ws.Range("D8").Value = "'0.8496"
$ws.Range("E8").Value = "'-0.78%"
$ws.Range("D9").Value = "'0.8598"
$ws.Range("E9").Value = "'-2.42%"
$ws.Range("D10").Value = "'0.1382"
$ws.Range("E10").Value = "'1.17%"
$ws.Range("D11").Value = "'0.07085"
$ws.Range("E11").Value = "'-0.82%"
$ws.Range("D12").Value = "'0.03235"
$ws.Range("E12").Value = "'12.90%"
$ws.Range("D13").Value = "'0.09361"
$ws.Range("E13").Value = "'-0.34%"
$ws.Range("D14").Value = "'0.001530"
$ws.Range("E14").Value = "'0.63%"
$ws.Range("B15").Value = "One"
$ws.Range("C15").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D15").Value = "'0.0005939"
$ws.Range("E15").Value = "'-1.34%"
$ws.Range("B16").Value = "TigerCash"
$ws.Range("C16").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D16").Value = "'0.005899"
$ws.Range("E16").Value = "'-2.12%"
$ws.Range("B17").Value = "LEO"
$ws.Range("C17").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D17").Value = "'3.546"
$ws.Range("E17").Value = "'1.62%"
$ws.Range("B18").Value = "BTSEToken"
$ws.Range("C18").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D18").Value = "'2.217"
$ws.Range("E18").Value = "'1.76%"
$ws.Range("B19").Value = "BitpandaEcosystemToken"
$ws.Range("C19").Value = "https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best"
$ws.Range("D19").Value = "'0.3144"
$ws.Range("E19").Value = "'-2.47%"
$ws.Range("B20").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C20").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D20").Value = "'0.03422"
$ws.Range("E20").Value = "'5.36%"
$ws.Range("B21").Value = "ProBitToken"
$ws.Range("C21").Value = "https://coinranking.com/coin/lQP4d6T2+probittoken-prob"
$ws.Range("D21").Value = "'0.1316"
$ws.Range("E21").Value = "'1.21%"
$ws.Range("B22").Value = "MCDex"
$ws.Range("C22").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
$ws.Range("D22").Value = "'3.493"
$ws.Range("E22").Value = "'-0.94%"
$ws.Range("B23").Value = "ZBToken"
$ws.Range("C23").Value = "https://coinranking.com/coin/CxmvOsCyENPso+zbtoken-zb"
$ws.Range("D23").Value = "'0.1410"
$ws.Range("E23").Value = "'2.16%"
$ws.Range("B24").Value = "CoinExToken"
$ws.Range("C24").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
$ws.Range("D24").Value = "'0.04109"
$ws.Range("E24").Value = "'-1.02%"
$ws.Range("D25").Value = "'0.001221"
$ws.Range("E25").Value = "'0.54%"
$ws.Range("D26").Value = "'0.004159"
$ws.Range("E26").Value = "'-7.44%"
$ws.Range("E27").Value = "'-0.88%"
$ws.Range("E28").Value = "'4.73%"
$ws.Range("D40").Value = "'0.03753"
$ws.Range("E40").Value = "'-0.81%"
$ws.Range("D41").Value = "'0.1070"
$ws.Range("E41").Value = "'-0.11%"
$ws.Range("E42").Value = "'-15.43%"
$ws.Range("E43").Value = "'-15.90%"
$ws.Range("D44").Value = "'0.01020"
$ws.Range("E44").Value = "'1.81%"
$ws.Range("D45").Value = "'0.00005481"
$ws.Range("E45").Value = "'7.69%"
$ws.Range("E46").Value = "'-0.01%"
$ws.Range("D47").Value = "'0.07099"
$ws.Range("E47").Value = "'-11.26%"
$ws.Range("D48").Value = "'0.002468"
$ws.Range("E48").Value = "'-10.70%"
$ws.Range("E49").Value = "'-0.01%"
$ws.Range("E50").Value = "'-0.01%"
